# Update IPC PO (C), DELTA (D) and DELTA^2 (E) columns for rows 2-51
# plus the TOTAL (row 52) and MSE (row 53) summary cells, reflecting
# the refactored weight-handling predictions (sliding_window_results_window_8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$poValues = @{
    2 = 29.88066855575936
    3 = 29.77401489754656
    4 = 29.64367529096509
    5 = 29.9770110492835
    6 = 30.15190796648162
    7 = 30.34663417843318
    8 = 30.62369824705883
    9 = 30.31539290832943
    10 = 30.48004311243192
    11 = 30.62951007234524
    12 = 30.69147273849411
    13 = 30.33463104955399
    14 = 30.37174292734279
    15 = 30.76360678343317
    16 = 31.06347961641817
    17 = 31.37022274551597
    18 = 31.36086790810678
    19 = 31.69297716748968
    20 = 32.00502054386315
    21 = 31.95227342423743
    22 = 31.60719355010644
    23 = 32.13576646462457
    24 = 32.14611717632918
    25 = 33.19851626968185
    26 = 32.87005672226958
    27 = 33.09926301179274
    28 = 33.11066246943691
    29 = 33.73068155870988
    30 = 34.50721488327389
    31 = 34.51580185188255
    32 = 35.03799556991351
    33 = 34.87663432728867
    34 = 35.17898363211293
    35 = 35.86378850647009
    36 = 36.26778913958069
    37 = 37.30114710264916
    38 = 37.8392741502586
    39 = 38.43957258246311
    40 = 39.34244782130198
    41 = 40.07310150521455
    42 = 40.34373512892738
    43 = 40.02935357895179
    44 = 41.19614329238238
    45 = 41.43468154294521
    46 = 41.23756675007412
    47 = 41.01664378554329
    48 = 41.4105672201997
    49 = 41.92653814877664
    50 = 42.23666827446134
    51 = 43.37985615845385
}

$deltaValues = @{
    2 = 0.04066855575935691
    3 = -0.03598510245344144
    4 = -0.2763247090349132
    5 = -0.002988950716499517
    6 = 0.1119079664816169
    7 = 0.1366341784331802
    8 = 0.4036982470588342
    9 = -0.06460709167056677
    10 = 0.04004311243192049
    11 = 0.1495100723452403
    12 = 0.001472738494104675
    13 = -0.415368950446009
    14 = -0.5682570726572109
    15 = -0.1863932165668274
    16 = 0.04347961641817122
    17 = 0.2502227455159662
    18 = 0.08086790810678224
    19 = 0.3129771674896844
    20 = 0.4250205438631482
    21 = 0.302273424237427
    22 = -0.2728064498935616
    23 = -0.1442335353754345
    24 = -0.3038828236708184
    25 = 0.3485162696818449
    26 = -0.02994327773041761
    27 = -0.000736988207265199
    28 = -0.2893375305630883
    29 = 0.03068155870987965
    30 = 0.4072148832738876
    31 = 0.1158018518825514
    32 = 0.1379955699135138
    33 = -0.4233656727113271
    34 = -0.5210163678870714
    35 = -0.4362114935299033
    36 = -0.5322108604193048
    37 = 0.001147102649163401
    38 = -0.06072584974140227
    39 = -0.06042741753689285
    40 = 0.4424478213019825
    41 = 0.6731015052145466
    42 = 0.4437351289273792
    43 = -0.07064642104820962
    44 = 0.5961432923823793
    45 = 0.534681542945215
    46 = 0.03756675007411303
    47 = -0.4833562144567125
    48 = -0.3894327798002948
    49 = -0.2734618512233666
    50 = -0.4633317255386586
    51 = -0.3201438415461553
}

$deltaSqValues = @{
    2 = 0.001653931427551922
    3 = 0.001294927598584677
    4 = 0.07635534482322946
    5 = [double]"8.933826385662973e-06"
    6 = 0.0125233929620507
    7 = 0.01866889871611011
    8 = 0.1629722746783755
    9 = 0.004174076294129019
    10 = 0.001603450853235425
    11 = 0.02235326173267898
    12 = [double]"2.168958672017706e-06"
    13 = 0.172531364994619
    14 = 0.3229161006249427
    15 = 0.03474243118212822
    16 = 0.001890477043871304
    17 = 0.062611422373548
    18 = 0.006539618561566977
    19 = 0.09795470736986596
    20 = 0.1806424627057263
    21 = 0.0913692230002195
    22 = 0.07442335910352835
    23 = 0.02080331272689671
    24 = 0.09234477052214972
    25 = 0.1214635902329484
    26 = 0.0008965998812409229
    27 = [double]"5.431516176479719e-07"
    28 = 0.08371620659234606
    29 = 0.000941358044867792
    30 = 0.1658239611597659
    31 = 0.01341006889942838
    32 = 0.01904277731575548
    33 = 0.1792384928303145
    34 = 0.2714580556062361
    35 = 0.1902804670875889
    36 = 0.2832483999482567
    37 = [double]"1.315844487717692e-06"
    38 = 0.003687628826815366
    39 = 0.003651472790177985
    40 = 0.195760074574871
    41 = 0.4530656363220883
    42 = 0.1969008646441978
    43 = 0.004990916806920915
    44 = 0.355386825052503
    45 = 0.2858843523662757
    46 = 0.001411260711130871
    47 = 0.2336332300539234
    48 = 0.1516578899829849
    49 = 0.07478138407451071
    50 = 0.2146762878906309
    51 = 0.1024920792799298
}

foreach ($row in 2..51) {
    $ws.Cells.Item($row, 3).Value = $poValues[$row]
    $ws.Cells.Item($row, 4).Value = $deltaValues[$row]
    $ws.Cells.Item($row, 5).Value = $deltaSqValues[$row]
}

# TOTAL row: sum of DELTA and sum of DELTA^2
$ws.Cells.Item(52, 3).Value = -0.5573866408334638
$ws.Cells.Item(52, 5).Value = 5.067881652051883

# MSE row: average of DELTA^2
$ws.Cells.Item(53, 5).Value = 0.1013576330410376

Write-Host "Updated IPC PO, DELTA, DELTA^2 columns and summary rows."
